$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.317.76'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.919.21'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("D5").Value = '0.8098'
$ws.Range("E5").Value = '  +3.54%  '
$ws.Range("D6").Value = '244.38'
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.3241'
$ws.Range("E8").Value = '  +2.74%  '
$ws.Range("D9").Value = '27.14'
$ws.Range("E9").Value = '  +3.90%  '
$ws.Range("D10").Value = '0.07102'
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("D11").Value = '0.7861'
$ws.Range("E11").Value = '  +6.17%  '
$ws.Range("D12").Value = '0.08093'
$ws.Range("E12").Value = '  +1.67%  '
$ws.Range("D13").Value = '1.903.63'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").Value = '5.419'
$ws.Range("E14").Value = '  +4.24%  '
$ws.Range("D15").Value = '94.77'
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").Value = '30.311.05'
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("D17").Value = '14.30'
$ws.Range("E17").Value = '  +2.89%  '
$ws.Range("D18").Value = '6.033'
$ws.Range("E18").Value = '  +2.94%  '
$ws.Range("D19").Value = '249.96'
$ws.Range("E19").Value = '  +1.94%  '
$ws.Range("D20").Value = '0.000007812'
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").Value = '2.171.52'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '7.983'
$ws.Range("E23").Value = '  +16.56%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '0.1618'
$ws.Range("E25").Value = '  +17.32%  '
$ws.Range("D26").Value = '9.504'
$ws.Range("E26").Value = '  +2.68%  '
$ws.Range("D27").Value = '167.65'
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("D28").Value = '19.10'
$ws.Range("E28").Value = '  +1.25%  '
$ws.Range("D29").Value = '2.134'
$ws.Range("E29").Value = '  +5.15%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '1.537'
$ws.Range("E31").Value = '  +1.29%  '
$ws.Range("D32").Value = '4.356'
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").Value = '0.05623'
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("D34").Value = '4.136'
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").Value = '1.302'
$ws.Range("E35").Value = '  +3.91%  '
$ws.Range("D36").Value = '0.7442'
$ws.Range("E36").Value = '  +1.71%  '
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("D38").Value = '2.717'
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").Value = '0.4484'
$ws.Range("E41").Value = '  +1.73%  '
$ws.Range("D42").Value = '73.74'
$ws.Range("E42").Value = '  +2.75%  '
$ws.Range("D43").Value = '5.977'
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("D44").Value = '0.8563'
$ws.Range("E44").Value = '  +1.85%  '
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '1.037.17'
$ws.Range("E47").Value = '  +5.81%  '
$ws.Range("D48").Value = '103.11'
$ws.Range("E48").Value = '  +2.77%  '
$ws.Range("D49").Value = '9.967'
$ws.Range("E49").Value = '  +2.42%  '
$ws.Range("D50").Value = '7.644'
$ws.Range("E50").Value = '  +1.72%  '
$ws.Range("D51").Value = '2.069.30'
$ws.Range("E51").Value = '  +0.72%  '
